# Fruta / hortaliza, semanal
# Insert a new weekly group of 3 rows (Especial / Primera / Segunda) at the
# top of the "Frutilla" data block (rows 448-450), pushing the previously
# existing 27 rows (9 weekly groups, rows 448-474) down by 3 rows
# (to rows 451-477). The new group uses the same static descriptive
# columns (A,B,C,E,F,G,H,I,J,K,Q,R,T) as the existing template rows, with
# an updated date and new Volumen/Precio/Precio-$-Kg figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 448:474 down to 451:477, inserting 3 blank rows at 448.
$ws.Rows("448:450").Insert()

$newRows = @(
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44516, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Especial",  440, 12500, 13000, 12750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1821, 7),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44516, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Primera",  360, 10500, 11000, 10750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1536, 7),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44516, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Segunda",  280,  8500,  9000,  8750, "`$/bandeja 7 kilos", "Provincia de Melipilla", 1250, 7)
)

$startRow = 448
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $rowVals[$c - 1]
    }
}
